$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 328
$ws.Range("F3").Value = 3323
$ws.Range("F6").Value = 7754
$ws.Range("F9").Value = 1149
$ws.Range("F10").Value = 1063
$ws.Range("F12").Value = 24
$ws.Range("F14").Value = 1766
$ws.Range("F16").Value = 99
$ws.Range("F17").Value = 2375
$ws.Range("F20").Value = 1032
$ws.Range("F22").Value = 6361
$ws.Range("F23").Value = 7005
$ws.Range("F24").Value = 407
$ws.Range("F26").Value = 1090
$ws.Range("F30").Value = 1079
$ws.Range("F32").Value = 521
$ws.Range("F33").Value = 521
$ws.Range("F38").Value = 603
$ws.Range("F41").Value = 1265
$ws.Range("F42").Value = 3265
$ws.Range("F45").Value = 481
$ws.Range("F49").Value = 478

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 648
$ws.Range("F10").Value = 295
$ws.Range("F25").Value = 9
$ws.Range("F26").Value = 6618
$ws.Range("F27").Value = 1
$ws.Range("F31").Value = 55

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2007
$ws.Range("F5").Value = 1326
$ws.Range("F8").Value = 2153
$ws.Range("F9").Value = 8943
$ws.Range("F10").Value = 1081
$ws.Range("F11").Value = 87

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3323
$ws.Range("F3").Value = 2007
$ws.Range("F5").Value = 7754
$ws.Range("F6").Value = 1326
$ws.Range("F8").Value = 1081
$ws.Range("F9").Value = 87
$ws.Range("F10").Value = 1149
$ws.Range("F11").Value = 1063
$ws.Range("F17").Value = 99
$ws.Range("F18").Value = 2375
$ws.Range("F20").Value = 1032
$ws.Range("F22").Value = 6361
$ws.Range("F23").Value = 7005
$ws.Range("F24").Value = 407
$ws.Range("F26").Value = 1090
$ws.Range("F31").Value = 521
$ws.Range("F36").Value = 603
$ws.Range("F40").Value = 3265
$ws.Range("F43").Value = 481
